$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.83567917099368
$ws.Range("C2").Value = 4.347754598935162
$ws.Range("D2").Value = 8.947934791546878
$ws.Range("E2").Value = 16.26076764462834
$ws.Range("F2").Value = 48.67805088739608
$ws.Range("K2").Value = 11.26397376326645

$ws.Range("B3").Value = 11.67813693589336
$ws.Range("C3").Value = 4.214988554051004
$ws.Range("D3").Value = 8.780091474645269
$ws.Range("E3").Value = 15.35067858270834
$ws.Range("F3").Value = 47.30757573533351
$ws.Range("K3").Value = 11.19011341761207

$ws.Range("B4").Value = 11.58636506302972
$ws.Range("C4").Value = 4.135341705279166
$ws.Range("D4").Value = 8.67456169930343
$ws.Range("E4").Value = 14.76924474061031
$ws.Range("F4").Value = 46.44854087176736
$ws.Range("K4").Value = 11.15011589207205

$ws.Range("B5").Value = 11.55026953716738
$ws.Range("C5").Value = 4.103436937274511
$ws.Range("D5").Value = 8.630961405431444
$ws.Range("E5").Value = 14.52689045055297
$ws.Range("F5").Value = 46.09446978730535
$ws.Range("K5").Value = 11.13517871219945

$ws.Range("B6").Value = 11.54435609469561
$ws.Range("C6").Value = 4.098174906308264
$ws.Range("D6").Value = 8.623686385271904
$ws.Range("E6").Value = 14.48632975882486
$ws.Range("F6").Value = 46.03544694541981
$ws.Range("K6").Value = 11.13278108349184

$ws.Range("B7").Value = 11.58587292335369
$ws.Range("C7").Value = 4.134909083424292
$ws.Range("D7").Value = 8.673976068866404
$ws.Range("E7").Value = 14.76599778755941
$ws.Range("F7").Value = 46.44378141327907
$ws.Range("K7").Value = 11.14990891000382

$ws.Range("B8").Value = 11.78036335212133
$ws.Range("C8").Value = 4.301636603088933
$ws.Range("D8").Value = 8.890590948904729
$ws.Range("E8").Value = 15.95181725065424
$ws.Range("F8").Value = 48.20939132980891
$ws.Range("K8").Value = 11.23740628242705

$ws.Range("B9").Value = 12.19852778953494
$ws.Range("C9").Value = 4.639975696130906
$ws.Range("D9").Value = 9.294669016677213
$ws.Range("E9").Value = 18.09520718372749
$ws.Range("F9").Value = 51.5159290358671
$ws.Range("K9").Value = 11.45060662414416

$ws.Range("B10").Value = 12.52453779673058
$ws.Range("C10").Value = 4.956420271253551
$ws.Range("D10").Value = 9.577567885519395
$ws.Range("E10").Value = 19.71836995631918
$ws.Range("F10").Value = 53.83016464184879
$ws.Range("K10").Value = 11.63127218317549

$ws.Range("B11").Value = 12.67612780060961
$ws.Range("C11").Value = 5.135137720577665
$ws.Range("D11").Value = 9.702968576233385
$ws.Range("E11").Value = 20.4159767303297
$ws.Range("F11").Value = 54.85430601257132
$ws.Range("K11").Value = 11.71834083024951

$ws.Range("B12").Value = 12.73393651624347
$ws.Range("C12").Value = 5.201349327230089
$ws.Range("D12").Value = 9.749960865130955
$ws.Range("E12").Value = 20.67431822616802
$ws.Range("F12").Value = 55.23773712079549
$ws.Range("K12").Value = 11.75198229709522

$ws.Range("B13").Value = 12.72146943886117
$ws.Range("C13").Value = 5.187154467467934
$ws.Range("D13").Value = 9.739862534991182
$ws.Range("E13").Value = 20.61893833800336
$ws.Range("F13").Value = 55.15535765127834
$ws.Range("K13").Value = 11.74470770854409

$ws.Range("B14").Value = 12.68087602963648
$ws.Range("C14").Value = 5.140614291093826
$ws.Range("D14").Value = 9.706844680352742
$ws.Range("E14").Value = 20.43734721311222
$ws.Range("F14").Value = 54.88594055689408
$ws.Range("K14").Value = 11.72109528440824

$ws.Range("B15").Value = 12.65606212215681
$ws.Range("C15").Value = 5.111916613627445
$ws.Range("D15").Value = 9.68655532929994
$ws.Range("E15").Value = 20.32535964171178
$ws.Range("F15").Value = 54.72033555600652
$ws.Range("K15").Value = 11.70671833722702

$ws.Range("B16").Value = 12.51469136335729
$ws.Range("C16").Value = 4.944535161009044
$ws.Range("D16").Value = 9.569304604059596
$ws.Range("E16").Value = 19.67196230152317
$ws.Range("F16").Value = 53.7626347222059
$ws.Range("K16").Value = 11.62567762878983

$ws.Range("B17").Value = 12.42875837701742
$ws.Range("C17").Value = 4.839232722154044
$ws.Range("D17").Value = 9.496516754917286
$ws.Range("E17").Value = 19.26070616827286
$ws.Range("F17").Value = 53.1675834113715
$ws.Range("K17").Value = 11.57719041949071

$ws.Range("B18").Value = 12.37964536786066
$ws.Range("C18").Value = 4.780592728263422
$ws.Range("D18").Value = 9.454342471934888
$ws.Range("E18").Value = 19.02032209695242
$ws.Range("F18").Value = 52.82264678763372
$ws.Range("K18").Value = 11.54976348566101

$ws.Range("B19").Value = 12.36307244055617
$ws.Range("C19").Value = 4.767846018382396
$ws.Range("D19").Value = 9.440010635262967
$ws.Range("E19").Value = 18.93827146204531
$ws.Range("F19").Value = 52.7054058739826
$ws.Range("K19").Value = 11.54055741972845

$ws.Range("B20").Value = 12.43787414974016
$ws.Range("C20").Value = 4.850542222961895
$ws.Range("D20").Value = 9.504297244345784
$ws.Range("E20").Value = 19.30488228760288
$ws.Range("F20").Value = 53.23120696085179
$ws.Range("K20").Value = 11.58230441343517

$ws.Range("B21").Value = 12.69278884607328
$ws.Range("C21").Value = 5.154323946790662
$ws.Range("D21").Value = 9.716556398066107
$ws.Range("E21").Value = 20.49084278024592
$ws.Range("F21").Value = 54.96519605525674
$ws.Range("K21").Value = 11.72801289047385

$ws.Range("B22").Value = 12.86171943648843
$ws.Range("C22").Value = 5.344329282107878
$ws.Range("D22").Value = 9.852390507646284
$ws.Range("E22").Value = 21.23198899797757
$ws.Range("F22").Value = 56.07276731222953
$ws.Range("K22").Value = 11.82713396140711

$ws.Range("B23").Value = 12.77136705321076
$ws.Range("C23").Value = 5.243697497079942
$ws.Range("D23").Value = 9.78016424548623
$ws.Range("E23").Value = 20.8395193924209
$ws.Range("F23").Value = 55.4840700462631
$ws.Range("K23").Value = 11.77388593053766

$ws.Range("B24").Value = 12.43375199731872
$ws.Range("C24").Value = 4.845432291676977
$ws.Range("D24").Value = 9.500780703808287
$ws.Range("E24").Value = 19.28492256224133
$ws.Range("F24").Value = 53.20245155938681
$ws.Range("K24").Value = 11.57999097582558

$ws.Range("B25").Value = 12.08183514572709
$ws.Range("C25").Value = 4.547697622136287
$ws.Range("D25").Value = 9.187716290591796
$ws.Range("E25").Value = 17.53147427128376
$ws.Range("F25").Value = 50.64028923727977
$ws.Range("K25").Value = 11.38860331065448

